$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new place row (row 13) for "Appartments" in Arabkir
$ws.Range("A13").Value = 12
$ws.Range("B13").Value = "Appartments"
$ws.Range("C13").Value = "Apartment"
$ws.Range("D13").Value = 15
$ws.Range("E13").Value = 40.2046591500417
$ws.Range("F13").Value = 44.5204283363223
$ws.Range("G13").Value = "Arabkir"

$ws.Range("C13").Select()
